# Export fixture update:
#  - the old "Stock" sheet becomes "Ubicaciones" (new warehouse-locations
#    export) and keeps its position (5th tab).
#  - a brand-new "Stock" sheet is appended right after it (6th tab), holding
#    the stock-by-location data that used to live directly on the old sheet.
#  - "Variantes" gains a first data row.

$wb = $excel.ActiveWorkbook

# --- 1. Sheet restructuring ------------------------------------------------
# Rename the existing "Stock" worksheet to "Ubicaciones" (same tab slot,
# same underlying sheet).
$ubicaciones = $wb.Worksheets.Item("Stock")
$ubicaciones.Name = "Ubicaciones"

# Insert a fresh worksheet right after it and call it "Stock".
$stock = $wb.Worksheets.Add($null, $ubicaciones)
$stock.Name = "Stock"

# --- 2. Variantes: add the first product-variant row -----------------------
$variantes = $wb.Worksheets.Item("Variantes")

$variantes.Cells.Item(2, 1).Value = 2
$variantes.Cells.Item(2, 2).Value = 1
$variantes.Cells.Item(2, 3).Value = "Nueces, estoy a dieta"
$variantes.Cells.Item(2, 4).Value = 12345678
$variantes.Cells.Item(2, 5).Value = 2000.0
$variantes.Cells.Item(2, 6).Value = 200.0
$variantes.Cells.Item(2, 7).Value = 10.0
$variantes.Cells.Item(2, 8).Value = 10.0
$variantes.Cells.Item(2, 9).Value = 10.0

# Widen columns C-G now that they hold real data (values chosen so the
# COM ColumnWidth->stored-width conversion lands on the target widths).
$variantes.Columns.Item(3).ColumnWidth = 16.68
$variantes.Columns.Item(4).ColumnWidth = 11.18
$variantes.Columns.Item(5).ColumnWidth = 7.88
$variantes.Columns.Item(6).ColumnWidth = 6.78
$variantes.Columns.Item(7).ColumnWidth = 5.68

# --- 3. Ubicaciones: new header row + first location row -------------------
$ubicaciones.Cells.Item(1, 1).Value = "ID"
$ubicaciones.Cells.Item(1, 2).Value = "Nombre"
$ubicaciones.Cells.Item(1, 3).Value = "Nombre Interno"
$ubicaciones.Cells.Item(1, 4).Value = "Calle"
$ubicaciones.Cells.Item(1, 5).Value = "Ciudad"
$ubicaciones.Cells.Item(1, 6).Value = "Calle de referencia"
$ubicaciones.Cells.Item(1, 7).Value = "Código Postal"
$ubicaciones.Cells.Item(1, 8).Value = "Teléfono"
$ubicaciones.Cells.Item(1, 9).Value = "País"
$ubicaciones.Cells.Item(1, 10).Value = "Región"
$ubicaciones.Cells.Item(1, 11).Value = "Activa"
$ubicaciones.Cells.Item(1, 12).Value = "Por defecto"
$ubicaciones.Cells.Item(1, 13).Value = "Backorderable"
$ubicaciones.Cells.Item(1, 14).Value = "Propagar por todas las variantes"

$ubicaciones.Cells.Item(2, 1).Value = 1
$ubicaciones.Cells.Item(2, 2).Value = "Isla Diamante"
$ubicaciones.Cells.Item(2, 3).Value = "Central"
$ubicaciones.Cells.Item(2, 4).Value = "Playa 123"
$ubicaciones.Cells.Item(2, 5).Value = "Til Til"
$ubicaciones.Cells.Item(2, 6).Value = "Juan algo 234"
$ubicaciones.Cells.Item(2, 7).Value = 12345
$ubicaciones.Cells.Item(2, 8).Value = 76543469
$ubicaciones.Cells.Item(2, 9).Value = "Chile"
$ubicaciones.Cells.Item(2, 10).Value = "Región Metropolitana"
$ubicaciones.Cells.Item(2, 11).Value = "Sí"
$ubicaciones.Cells.Item(2, 12).Value = "No"
$ubicaciones.Cells.Item(2, 13).Value = "Sí"
$ubicaciones.Cells.Item(2, 14).Value = "Sí"

$ubicaciones.Columns.Item(1).ColumnWidth = 3.48
$ubicaciones.Columns.Item(2).ColumnWidth = 7.88
$ubicaciones.Columns.Item(3).ColumnWidth = 13.38
$ubicaciones.Columns.Item(4).ColumnWidth = 8.98
$ubicaciones.Columns.Item(5).ColumnWidth = 6.78
$ubicaciones.Columns.Item(6).ColumnWidth = 13.38
$ubicaciones.Columns.Item(7).ColumnWidth = 11.18
$ubicaciones.Columns.Item(8).ColumnWidth = 11.18
$ubicaciones.Columns.Item(9).ColumnWidth = 5.68
$ubicaciones.Columns.Item(10).ColumnWidth = 15.58
$ubicaciones.Columns.Item(11).ColumnWidth = 4.58
$ubicaciones.Columns.Item(12).ColumnWidth = 10.08
$ubicaciones.Columns.Item(13).ColumnWidth = 8.98
$ubicaciones.Columns.Item(14).ColumnWidth = 18.88

# --- 4. Stock (new sheet): header row + first stock-by-location row --------
$stock.Cells.Item(1, 1).Value = "ID Producto"
$stock.Cells.Item(1, 2).Value = "ID Variante"
$stock.Cells.Item(1, 3).Value = "Cantidad"
$stock.Cells.Item(1, 4).Value = "ID Ubicación"
$stock.Cells.Item(1, 5).Value = "Backorderable"

$stock.Cells.Item(2, 1).Value = 1
$stock.Cells.Item(2, 2).Value = 2
$stock.Cells.Item(2, 3).Value = 20
$stock.Cells.Item(2, 4).Value = 1
$stock.Cells.Item(2, 5).Value = "Sí"

$stock.Columns.Item(1).ColumnWidth = 10.08
$stock.Columns.Item(2).ColumnWidth = 7.88
$stock.Columns.Item(3).ColumnWidth = 6.78
$stock.Columns.Item(4).ColumnWidth = 8.98
$stock.Columns.Item(5).ColumnWidth = 8.98
